$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the first paragraph entirely ("Git" + " Test of the doc file",
#    wrapped with a spellcheck proofErr run around "Git"). Its own paragraph
#    mark is removed along with it, so the document collapses back down to
#    just the "中华人民共和国" paragraph.
# ---------------------------------------------------------------------------
$first = $d.Paragraphs(1)
$first.Range.Delete()

# ---------------------------------------------------------------------------
# 2. The remaining paragraph ("中华人民共和国", with the eastAsia rFonts hint
#    on its pPr and the _GoBack bookmark) needs to turn into four paragraphs:
#      a) "Git Test" (plain run, no run-level rPr) + the _GoBack bookmark,
#         reusing the original pPr/rFonts hint.
#      b) "中" and "国" as two separate eastAsia-hinted runs.
#      c) an empty paragraph.
#      d) "The change happen" + " in the word" + "！" as three eastAsia-hinted
#         runs.
#    We build all four as fresh OOXML and insert them together right after
#    the still-present old paragraph (targeting a position strictly inside
#    its text keeps the insertion anchored "after this paragraph" instead of
#    merging into it), then delete the old paragraph afterward. That gives
#    every new run a clean slate rPr instead of inheriting the old run's
#    eastAsia hint.
# ---------------------------------------------------------------------------
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$para1 = "<w:p $w>" +
           "<w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr>" +
           "<w:r><w:t>Git Test</w:t></w:r>" +
           "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
           "<w:bookmarkEnd w:id=`"0`"/>" +
         "</w:p>"

$para2 = "<w:p $w>" +
           "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>中</w:t></w:r>" +
           "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>国</w:t></w:r>" +
         "</w:p>"

$para3 = "<w:p $w/>"

$para4 = "<w:p $w>" +
           "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>The change happen</w:t></w:r>" +
           "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t xml:space=`"preserve`"> in the word</w:t></w:r>" +
           "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>！</w:t></w:r>" +
         "</w:p>"

$batch = $para1 + $para2 + $para3 + $para4

$old = $d.Paragraphs(1)
$anchorPos = $old.Range.Start + 1
$anchor = $d.Range($anchorPos, $anchorPos)
$anchor.InsertXML($batch)

# Now remove the original "中华人民共和国" paragraph (still paragraph 1).
$old2 = $d.Paragraphs(1)
$old2.Range.Delete()

"done"
